$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 220 (shifts existing rows 220:288 down to 221:289)
$ws.Rows("220").Insert()

# Populate the new row 220 with the new weekly record
$ws.Range("A220").Value = 3
$ws.Range("B220").Value = 'Femacal de La Calera'
$ws.Range("C220").Value = 'Coquimbo'
$ws.Range("D220").Value = 44588
$ws.Range("E220").Value = 5
$ws.Range("F220").Value = 100112040
$ws.Range("G220").Value = 'Cilantro'
$ws.Range("H220").Value = 'Sin especificar'
$ws.Range("I220").Value = 'Primera'
$ws.Range("J220").Value = 230
$ws.Range("K220").Value = 4000
$ws.Range("L220").Value = 4300
$ws.Range("M220").Value = 4143
$ws.Range("N220").Value = '$/docena de atados (3 kilos)'
$ws.Range("O220").Value = 'Provincia de Quillota'
$ws.Range("P220").Value = 1381
$ws.Range("Q220").Value = 3
$ws.Range("R220").Value = 'Hortaliza'
